# Nexial "base-showcase" workbook update:
#   - [JSON] add `storeKeys(json,jsonpath,var)` (extract immediate keys of a
#     resolved JSON fragment based on jsonpath) to the #system function list.
#   - the single-entry `text` category (`spellCheck(var,profile,text)`) is
#     retired, so its data column is removed and everything to its right
#     shifts left by one column; the `target` category list loses the
#     now-unused `text` row and everything below it shifts up by one row.
#
# NOTE: in this host, Range.Insert()/Range.Delete() operate on the whole
# row (all columns), not just the addressed column, so the column-A and
# column-M single-column shifts below are done by direct value assignment
# instead of Insert/Delete -- only the real column-wide operation
# (deleting the entire "Y" column) uses Columns(...).Delete().

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Drop the "text" function-category column (Y). Everything from "web"
#    (previously column Z) onward shifts one column to the left.
# ---------------------------------------------------------------------
$ws.Columns("Y:Y").Delete()

# ---------------------------------------------------------------------
# 2) Remove the "text" row from the category list in column A (row 25),
#    shifting the remaining categories (web, webalert, webcookie, ws,
#    ws.async, xml) up by one row; the list is now one row shorter.
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "web"
$ws.Range("A26").Value = "webalert"
$ws.Range("A27").Value = "webcookie"
$ws.Range("A28").Value = "ws"
$ws.Range("A29").Value = "ws.async"
$ws.Range("A30").Value = "xml"
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------
# 3) Insert the new `storeKeys(json,jsonpath,var)` function into the
#    alphabetically-sorted "json" list, between storeCount(...) and
#    storeValue(...) -- i.e. at M16, pushing the two storeValue* entries
#    down by one row.
# ---------------------------------------------------------------------
$ws.Range("M18").Value = "storeValues(json,jsonpath,var)"
$ws.Range("M17").Value = "storeValue(json,jsonpath,var)"
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 4) Fix up the defined names whose extents moved because of the above
#    structural edits (this host does not auto-adjust named ranges on
#    row/column delete, so it's done explicitly). "text" itself is left
#    untouched -- it's simply not cleaned up, still pointing at $Y$2:$Y$2.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
